$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4578.625
$ws.Range("I38").Value = 3725.3
$ws.Range("K38").Value = 11175.9
$ws.Range("M38").Value = -10803.9
$ws.Range("N40").ClearContents()
$ws.Range("H40").Value = 2666.3333
$ws.Range("I40").Value = 2666.3333
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2666.3333
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2491.3333
$ws.Range("H64").Value = 9904.454
$ws.Range("I64").Value = 9509.799999999999
$ws.Range("K64").Value = 9509.799999999999
$ws.Range("M64").Value = -9261.799999999999
$ws.Range("H67").Value = 9904.454
$ws.Range("I67").Value = 9509.799999999999
$ws.Range("K67").Value = 9509.799999999999
$ws.Range("M67").Value = -8651.799999999999
$ws.Range("H76").Value = 15536.546
$ws.Range("I76").Value = 24580.4
$ws.Range("K76").Value = 24580.4
$ws.Range("M76").Value = -24265.4
$ws.Range("H79").Value = 15536.546
$ws.Range("I79").Value = 24580.4
$ws.Range("K79").Value = 24580.4
$ws.Range("M79").Value = -23488.4
$ws.Range("N99").ClearContents()
$ws.Range("H99").Value = 702.125
$ws.Range("I99").Value = 702.125
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2106.375
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -608.375
$ws.Range("H100").Value = 1679.1666
$ws.Range("I100").Value = 901.8182
$ws.Range("K100").Value = 901.8182
$ws.Range("M100").Value = -360.8182
$ws.Range("H106").Value = 1806.6471
$ws.Range("I106").Value = 1514.2
$ws.Range("K106").Value = 1514.2
$ws.Range("M106").Value = -883.2
$ws.Range("H107").Value = 865.0476
$ws.Range("I107").Value = 601.125
$ws.Range("K107").Value = 601.125
$ws.Range("M107").Value = 1318.875
$ws.Range("N115").Value = -4634
$ws.Range("H115").Value = 344.6
$ws.Range("J115").Value = 500
$ws.Range("L115").Value = 1500
$ws.Range("H127").Value = 2526
$ws.Range("I127").Value = 1399
$ws.Range("K127").Value = 4197
$ws.Range("M127").Value = 763
$ws.Range("H132").Value = 1916.965
$ws.Range("I132").Value = 1566.9183
$ws.Range("K132").Value = 4700.7549
$ws.Range("M132").Value = -2170.7549
$ws.Range("H137").Value = 2495.3333
$ws.Range("I137").Value = 2617.087
$ws.Range("J137").Value = 2320.3125
$ws.Range("K137").Value = 7851.261
$ws.Range("L137").Value = 6960.9375
$ws.Range("M137").Value = -5301.261
$ws.Range("N137").Value = -12060.9375
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 43482910
$ws.Range("I61").Value = 58826852
$ws.Range("K61").Value = 58826852
$ws.Range("M61").Value = -58826640
$ws.Range("H122").Value = 2893.1667
$ws.Range("J122").Value = 4475.2
$ws.Range("L122").Value = 13425.6
$ws.Range("N122").Value = -18325.6
$ws.Range("H132").Value = 30306422
$ws.Range("I132").Value = 3482.1614
$ws.Range("J132").Value = 500002000
$ws.Range("K132").Value = 10446.4842
$ws.Range("L132").Value = 1500006000
$ws.Range("M132").Value = -7916.484199999999
$ws.Range("N132").Value = -1500011060
$ws.Range("H136").Value = 43482910
$ws.Range("I136").Value = 58826852
$ws.Range("K136").Value = 176480556
$ws.Range("M136").Value = -176478006
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2887.024
$ws.Range("I134").Value = 2752.718
$ws.Range("K134").Value = 8258.153999999999
$ws.Range("M134").Value = -5723.153999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 339.35
$ws.Range("J7").Value = 361.23077
$ws.Range("L7").Value = 361.23077
$ws.Range("N7").Value = -587.23077
$ws.Range("H107").Value = 1119.5
$ws.Range("I107").Value = 556.0909
$ws.Range("K107").Value = 556.0909
$ws.Range("M107").Value = 1363.9091
$ws.Range("H122").Value = 1842.174
$ws.Range("I122").Value = 1766.0588
$ws.Range("J122").Value = 2057.8333
$ws.Range("K122").Value = 5298.1764
$ws.Range("L122").Value = 6173.499899999999
$ws.Range("M122").Value = -2848.1764
$ws.Range("N122").Value = -11073.4999
$ws.Range("H132").Value = 3573.5833
$ws.Range("J132").Value = 5725.1665
$ws.Range("L132").Value = 17175.4995
$ws.Range("N132").Value = -22235.4995
$ws.Range("H134").Value = 1395.7646
$ws.Range("I134").Value = 1182.0667
$ws.Range("K134").Value = 3546.2001
$ws.Range("M134").Value = -1011.2001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2526.5
$ws.Range("I3").Value = 789.75
$ws.Range("K3").Value = 2369.25
$ws.Range("M3").Value = -2257.25
$ws.Range("H11").Value = 149.5
$ws.Range("I11").Value = 99.333336
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 298.000008
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = -158.000008
$ws.Range("N11").Value = -1180
$ws.Range("H131").Value = 1708
$ws.Range("I131").Value = 1287.5
$ws.Range("J131").Value = 1948.2858
$ws.Range("K131").Value = 3862.5
$ws.Range("L131").Value = 5844.857400000001
$ws.Range("M131").Value = 1177.5
$ws.Range("N131").Value = -15924.8574
$ws.Range("H132").Value = 3032899.8
$ws.Range("I132").Value = 1827.2858
$ws.Range("J132").Value = 4447400
$ws.Range("K132").Value = 16445.5722
$ws.Range("L132").Value = 40026600
$ws.Range("M132").Value = -13915.5722
$ws.Range("N132").Value = -40031660
$ws.Range("H139").Value = 2745.2666
$ws.Range("J139").Value = 3145.125
$ws.Range("L139").Value = 9435.375
$ws.Range("N139").Value = -19715.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8479262
$ws.Range("I11").Value = 9371105
$ws.Range("J11").Value = 6750
$ws.Range("K11").Value = 9371105
$ws.Range("L11").Value = 6750
$ws.Range("M11").Value = -9370966
$ws.Range("N11").Value = -7028
$ws.Range("H70").Value = 8501.4
$ws.Range("I70").Value = 4749
$ws.Range("J70").Value = 11003
$ws.Range("K70").Value = 4749
$ws.Range("L70").Value = 11003
$ws.Range("M70").Value = -4479
$ws.Range("N70").Value = -11543
$ws.Range("H73").Value = 8501.4
$ws.Range("I73").Value = 4749
$ws.Range("J73").Value = 11003
$ws.Range("K73").Value = 4749
$ws.Range("L73").Value = 11003
$ws.Range("M73").Value = -3813
$ws.Range("N73").Value = -12875
$ws.Range("H97").Value = 1312.8182
$ws.Range("I97").Value = 444.2
$ws.Range("K97").Value = 444.2
$ws.Range("M97").Value = 51.80000000000001
$ws.Range("H102").Value = 3719
$ws.Range("I102").Value = 2800
$ws.Range("J102").Value = 3902.8
$ws.Range("K102").Value = 2800
$ws.Range("L102").Value = 3902.8
$ws.Range("M102").Value = -1178
$ws.Range("N102").Value = -7146.8
$ws.Range("H113").Value = 3232.875
$ws.Range("I113").Value = 2324.1875
$ws.Range("J113").Value = 5050.25
$ws.Range("K113").Value = 2324.1875
$ws.Range("L113").Value = 5050.25
$ws.Range("M113").Value = -154.1875
$ws.Range("N113").Value = -9390.25
$ws.Range("H122").Value = 7609.294
$ws.Range("J122").Value = 9156.666999999999
$ws.Range("L122").Value = 27470.001
$ws.Range("N122").Value = -32370.001
$ws.Range("H132").Value = 3355.889
$ws.Range("I132").Value = 3084.8948
$ws.Range("K132").Value = 9254.6844
$ws.Range("M132").Value = -6724.6844
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3894.8708
$ws.Range("I7").Value = 3888.7
$ws.Range("J7").Value = 3906.0908
$ws.Range("K7").Value = 3888.7
$ws.Range("L7").Value = 3906.0908
$ws.Range("M7").Value = -3776.7
$ws.Range("N7").Value = -4130.0908
$ws.Range("N16").Value = -515
$ws.Range("H16").Value = 619.4286
$ws.Range("I16").Value = 797.2
$ws.Range("J16").Value = 175
$ws.Range("K16").Value = 797.2
$ws.Range("L16").Value = 175
$ws.Range("M16").Value = -627.2
$ws.Range("H22").Value = 2913
$ws.Range("I22").Value = 2010
$ws.Range("J22").Value = 3235.5
$ws.Range("K22").Value = 2010
$ws.Range("L22").Value = 3235.5
$ws.Range("M22").Value = -1715
$ws.Range("N22").Value = -3825.5
$ws.Range("H27").Value = 2913
$ws.Range("I27").Value = 2010
$ws.Range("J27").Value = 3235.5
$ws.Range("K27").Value = 2010
$ws.Range("L27").Value = 3235.5
$ws.Range("M27").Value = -1903
$ws.Range("N27").Value = -3449.5
$ws.Range("N43").Value = -28352
$ws.Range("H43").Value = 27966
$ws.Range("J43").Value = 27966
$ws.Range("L43").Value = 27966
$ws.Range("H46").Value = 1340.6531
$ws.Range("I46").Value = 664.8857400000001
$ws.Range("J46").Value = 3030.0715
$ws.Range("K46").Value = 664.8857400000001
$ws.Range("L46").Value = 3030.0715
$ws.Range("M46").Value = -476.8857400000001
$ws.Range("N46").Value = -3406.0715
$ws.Range("H126").Value = 3894.8708
$ws.Range("I126").Value = 3888.7
$ws.Range("J126").Value = 3906.0908
$ws.Range("K126").Value = 11666.1
$ws.Range("L126").Value = 11718.2724
$ws.Range("M126").Value = -9196.099999999999
$ws.Range("N126").Value = -16658.2724
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M30").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("H122").Value = 62502150
$ws.Range("I122").Value = 76924880
$ws.Range("K122").Value = 230774640
$ws.Range("M122").Value = -230772190
